$d = $word.ActiveDocument

# 1) "Entities: Nodes. Data. Resources. Models / State. ..." bullet:
#    "Object From Schema (Concept) to Behavior (Sign) Mappings."
# -> "Arcs Recognition (Behavior). Objects From Schema. (Concept) to Behavior (Sign) Mappings."
$d.Content.Find.Execute(
    "Object From Schema (Concept) to Behavior (Sign) Mappings.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Arcs Recognition (Behavior). Objects From Schema. (Concept) to Behavior (Sign) Mappings.",
    2
) | Out-Null

# 2) "Entities: Type / Context. Schema. Kinds / Roles. ..." bullet:
#    "Arcs Grammar / Recognition."
# -> "Node Recognition (Data)."
$d.Content.Find.Execute(
    "Arcs Grammar / Recognition.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Node Recognition (Data).",
    2
) | Out-Null

# 3) "Entities: Arcs. Behavior. Statements Kinds / NER Aggregation. ..." bullet:
#    "Context Type / Relationship / Dimension Recognition."
# -> "Context Types / Relationships / Dimensions Recognition (Schema)."
$d.Content.Find.Execute(
    "Context Type / Relationship / Dimension Recognition.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Context Types / Relationships / Dimensions Recognition (Schema).",
    2
) | Out-Null

# 4) The empty list paragraph immediately following the "Entities: Arcs. ..."
#    bullet carries an explicit "no underline" direct format on its paragraph
#    mark (<w:u w:val="none"/> inside w:pPr/w:rPr). The edit clears that
#    direct formatting. Locate the "Entities: Arcs..." paragraph by index
#    (Paragraph.Next is unreliable in this host, so use Paragraphs.Item),
#    then rewrite the following paragraph (same numbering/indent, no
#    explicit underline) via InsertXML - Font.Underline always re-serializes
#    an explicit value and can't express "no direct formatting".
$idx = 0
$foundIdx = -1
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -match "^Entities: Arcs\. Behavior\. Statements Kinds") {
        $foundIdx = $idx
        break
    }
}
if ($foundIdx -gt 0) {
    $target = $d.Paragraphs.Item($foundIdx + 1)
    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:pPr>' +
           '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>' +
           '<w:ind w:left="600" w:hanging="360"/>' +
           '</w:pPr>' +
           '<w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r>' +
           '</w:p>'
    $target.Range.InsertXML($xml)
}
